$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75 ---------------------------------------------------------
# Force column B to stay text (otherwise Excel auto-converts the
# "yyyy-mm-dd" looking string into a date serial number).
$ws.Cells.Item(75, 2).NumberFormat = "@"
$ws.Cells.Item(75, 2).Value = "2022-06-07"

# Copy the formatting of the last existing data row onto the new row so
# the new cells pick up the same cell style (s="2") used throughout the
# table, overriding the text number format applied above.
$ws.Range("A74:F74").Copy()
$ws.Range("A75:F75").PasteSpecial(-4122)

$ws.Cells.Item(75, 1).Value = 1654615201
$ws.Cells.Item(75, 3).Value = "通知 | 关于开展2022年上海交通大学暑期社会实践的通知"
$ws.Cells.Item(75, 4).Value = "JI青团"
$ws.Cells.Item(75, 5).Value = "https://mmbiz.qlogo.cn/mmbiz_png/QfDapvG9u4AmJ3iaZRWu8ZlB7ficzAEoXM0o3NLYG53goB7nr1kuOwvEicibDXWdD9Btq7Fwkj92KUMQDsjYsnd6qA/0?wx_fmt=png"
$ws.Cells.Item(75, 6).Value = "http://mp.weixin.qq.com/s?__biz=MzUyMzMyNTY0OQ==&mid=2247486349&idx=1&sn=4be0b005228f81c99e38b5610c9aa769&chksm=fa3f1233cd489b25f15a87cdaf3173ffef923601b0096c30b2e20bc6157e14f730649f861680#rd"

# --- Row 76 ---------------------------------------------------------
$ws.Cells.Item(76, 2).NumberFormat = "@"
$ws.Cells.Item(76, 2).Value = "2022-06-06"

$ws.Range("A74:F74").Copy()
$ws.Range("A76:F76").PasteSpecial(-4122)

$ws.Cells.Item(76, 1).Value = 1654524315
$ws.Cells.Item(76, 3).Value = "新一期青年大学习来啦！"
$ws.Cells.Item(76, 4).Value = "JI青团"
$ws.Cells.Item(76, 5).Value = "https://mmbiz.qlogo.cn/mmbiz_jpg/QfDapvG9u4DqIu8iccicyWEOzLmgsEGDsWKx0brvl21p54pndW57KibsUmWbgRYcibIRWa8b2xBmg6e6HRZx9rNaicA/0?wx_fmt=jpeg"
$ws.Cells.Item(76, 6).Value = "http://mp.weixin.qq.com/s?__biz=MzUyMzMyNTY0OQ==&mid=2247486337&idx=1&sn=37f520d21e21fc4cfa82a1056a7f9690&chksm=fa3f123fcd489b298963676544f00da8e408bfabdd7dedff9d7975ba67606f045fd0ee727048#rd"

# Refresh the sheet's used-range selection / dimension to span the two
# freshly appended rows, mirroring the original author's saved view.
$null = $ws.Range("A1:F76").Select()
